# Adds three new worksheets (CypherOutput_Message, StatOutput, StatOutput_Message)
# mirroring the existing CypherOutput / Message sheets, reporting file/sample/case/
# study counts returned by the companion "stat" Cypher query.

$wb = $excel.ActiveWorkbook

$neo4jUrlLbl  = 'Neo4j_URL:'
$neo4jUrlVal  = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$userNameLbl  = 'User_name:'
$userNameVal  = 'neo4j'
$pwdLbl       = 'PWD:'
$pwdVal       = 'icdcDBneo4j0'
$cypherLbl    = 'Cypher:'
$cypherOrig   = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Lymphoma :: Stage 2''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$cypherStat   = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Lymphoma :: Stage 2'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$outputLbl    = 'Output:'
$outputVal    = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC06_Canine_Filter_Diagnosis-LymphStg2_Neo4jData.xlsx'

function Write-MessageSheet($sheet, $cypherText) {
    $sheet.Range("A1").Value = $neo4jUrlLbl
    $sheet.Range("A2").Value = $neo4jUrlVal
    $sheet.Range("A3").Value = $userNameLbl
    $sheet.Range("A4").Value = $userNameVal
    $sheet.Range("A5").Value = $pwdLbl
    $sheet.Range("A6").Value = $pwdVal
    $sheet.Range("A7").Value = $cypherLbl
    $sheet.Range("A8").Value = $cypherText
    $sheet.Range("A9").Value = $outputLbl
    $sheet.Range("A10").Value = $outputVal
}

# --- Sheet 3: CypherOutput_Message ------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$cypherOutputMessage.Name = "CypherOutput_Message"
Write-MessageSheet $cypherOutputMessage $cypherOrig

# --- Sheet 4: StatOutput -----------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutput = $wb.Worksheets.Add($null, $lastSheet)
$statOutput.Name = "StatOutput"

$statOutput.Range("A1").Value = "number_of_files"
$statOutput.Range("B1").Value = "number_of_sample"
$statOutput.Range("C1").Value = "number_of_cases"
$statOutput.Range("D1").Value = "number_of_study"

# The counts look like plain numbers but must be stored as text (shared
# strings), matching the report generator's output. Writing them first as a
# formula and then converting the formula to its resulting value keeps the
# cell's type as text without forcing a new "quoted text" cell style.
$statValues = @("2", "5", "2", "1")
for ($col = 1; $col -le 4; $col++) {
    $cell = $statOutput.Cells.Item(2, $col)
    $cell.Formula = '="' + $statValues[$col - 1] + '"'
}
$statRow = $statOutput.Range("A2:D2")
$statRow.Copy()
$statRow.PasteSpecial(-4163) # xlPasteValues

# --- Sheet 5: StatOutput_Message ---------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$statOutputMessage.Name = "StatOutput_Message"
Write-MessageSheet $statOutputMessage $cypherStat

# Restore the original active sheet/tab selection (CypherOutput) so the
# workbook-level view state is left exactly as it was before the edit.
$wb.Worksheets.Item(1).Activate()

Write-Host "Added CypherOutput_Message, StatOutput, StatOutput_Message sheets"
